$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.878.92"
$ws.Range("E2").Value = "  +1.74%  "

$ws.Range("D3").Value = "1.728.21"
$ws.Range("E3").Value = "  +0.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9972"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.81"
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9977"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4886"
$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2598"
$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06220"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").Value = "1.732.39"
$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "15.99"
$ws.Range("E11").Value = "  +3.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06911"
$ws.Range("E12").Value = "  -1.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6085"
$ws.Range("E13").Value = "  +1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.488"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.18"
$ws.Range("E15").Value = "  -0.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9986"
$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").Value = "26.627.41"
$ws.Range("E17").Value = "  +0.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9973"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007178"
$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("D21").Value = "1.953.70"
$ws.Range("E21").Value = "  +0.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.426"
$ws.Range("E22").Value = "  -1.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.553"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.112"
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.52"
$ws.Range("E25").Value = "  +0.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.29"
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.778"
$ws.Range("E27").Value = "  +4.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.380"
$ws.Range("E28").Value = "  -1.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.11"
$ws.Range("E29").Value = "  -0.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.948"
$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07992"
$ws.Range("E31").Value = "  +0.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.685"
$ws.Range("E32").Value = "  +0.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04533"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9970"
$ws.Range("E34").Value = "  -0.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.596"
$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.008"
$ws.Range("E36").Value = "  +1.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6243"
$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9333"
$ws.Range("E38").Value = "  +1.83%  "

$ws.Range("E39").Value = "  +5.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.462"
$ws.Range("E40").Value = "  +2.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9987"
$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01502"
$ws.Range("E42").Value = "  +1.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.654"
$ws.Range("E43").Value = "  +5.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.43"
$ws.Range("E44").Value = "  -0.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3854"
$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.918"
$ws.Range("E46").Value = "  +3.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1160"
$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05389"
$ws.Range("E48").Value = "  +0.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.905"
$ws.Range("E49").Value = "  +2.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.14"
$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.238"
$ws.Range("E51").Value = "  -0.06%  "
